# Add a 4th adjacency-matrix table ("H-Graph") below the existing three,
# mirroring the layout/formatting of the table at A18:H23.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Tabelle1")

# 1) Clone the formatting (fills, borders, alignment, rotated header, merges'
#    base xf) of the third table onto the new block.
$ws.Range("A18:H23").Copy()
$ws.Range("A26").PasteSpecial(-4122)   # xlPasteFormats
$excel.CutCopyMode = 0

# The source table's rightmost two columns (G:H) carry an extra
# vertical-center tweak that the new table does not use - align G28:H31
# back to the same (non vertical-centered) look as C28:F31.
$ws.Range("G28:H31").VerticalAlignment = -4107   # xlBottom (= "no override")

# 2) Header labels (same shared strings as the other three tables).
$ws.Range("A26").Value = "N-Graph"
$ws.Range("B26").Value = "H-Graph"

# 3) Merge the label cells exactly like the other tables.
$ws.Range("A26:A31").Merge()
$ws.Range("B26:H26").Merge()

# 4) Column index header (1..6) for the new matrix.
$ws.Range("C27").Value = 1
$ws.Range("D27").Value = 2
$ws.Range("E27").Value = 3
$ws.Range("F27").Value = 4
$ws.Range("G27").Value = 5
$ws.Range("H27").Value = 6

# 5) Row index labels + adjacency matrix values.
$ws.Range("B28").Value = 1
$ws.Range("B29").Value = 2
$ws.Range("B30").Value = 3
$ws.Range("B31").Value = 4

$matrix = @(
  @(0, 1, 0, 0, 1, 0),
  @(0, 0, 1, 1, 0, 0),
  @(0, 0, 1, 1, 0, 0),
  @(0, 1, 0, 0, 1, 0)
)
for ($i = 0; $i -lt 4; $i++) {
  $row = 28 + $i
  for ($j = 0; $j -lt 6; $j++) {
    $col = 3 + $j
    $ws.Cells.Item($row, $col).Value = $matrix[$i][$j]
  }
}

# 6) Move the visible selection to where the author left off editing.
$ws.Range("E33").Select()
